# Insert a new data row before the existing row 586 (Vega Modelo de Temuco -
# Zapallo italiano price record), shifting rows 586:679 down to 587:680 and
# extending the used range from A1:R679 to A1:R680.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 586; everything below shifts down by one.
$ws.Rows.Item(586).Insert()

# Populate the newly inserted row 586 with the new price record.
$ws.Cells.Item(586, 1).Value  = 10
$ws.Cells.Item(586, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(586, 3).Value  = "La Araucanía"
$ws.Cells.Item(586, 4).Value  = 44951
$ws.Cells.Item(586, 5).Value  = 9
$ws.Cells.Item(586, 6).Value  = 100112032
$ws.Cells.Item(586, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(586, 8).Value  = "Sin especificar"
$ws.Cells.Item(586, 9).Value  = "Primera"
$ws.Cells.Item(586, 10).Value = 125
$ws.Cells.Item(586, 11).Value = 12000
$ws.Cells.Item(586, 12).Value = 12000
$ws.Cells.Item(586, 13).Value = 12000
$ws.Cells.Item(586, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(586, 15).Value = "Región del Maule"
$ws.Cells.Item(586, 16).Value = 240
$ws.Cells.Item(586, 17).Value = 50
$ws.Cells.Item(586, 18).Value = "Hortaliza"
